$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two kiosk rows for center 10001 (rows 3 and 4) used slightly different
# latitude/longitude values from row 2. Align them onto the same coordinates
# as row 2 (which also merges/dedupes the shared-string table, since the two
# stray longitude strings become unused).
$nbsp = [char]0xA0
$lon = "$nbsp-6.453275"

# Use a scratch cell with a text-valued formula + paste-special-values so the
# longitude string is written verbatim (t="s") instead of being silently
# re-parsed/auto-converted into a floating point number by plain value
# assignment.
$ws.Range("Z1").Formula = '="' + $lon + '"'
$ws.Range("Z1").Copy()
$ws.Range("H3").PasteSpecial(-4163)
$ws.Range("H4").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

$ws.Range("G3").Value = 34.521169999999998
$ws.Range("G4").Value = 34.521169999999998

# number_of_kiosks (column L): center 10001 now has 3 kiosks, every other
# center now has 2 kiosks.
$ws.Range("L2").Value = 3
$ws.Range("L3").Value = 3
$ws.Range("L4").Value = 3

for ($r = 5; $r -le 46; $r++) {
    $ws.Range("L$r").Value = 2
}
